# Update cryptocurrency price/volume data and reorder several coin rows
# per the Feb 13 2023 GitHub Actions data refresh.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'292.30"
$ws.Range("E2").Value = "'-5.47%"
$ws.Range("D3").Value = "'40.28"
$ws.Range("E3").Value = "'-1.98%"
$ws.Range("D4").Value = "'5.024"
$ws.Range("E4").Value = "'-3.68%"
$ws.Range("D5").Value = "'0.07393"
$ws.Range("E5").Value = "'-3.85%"
$ws.Range("B6").Value = "GateToken"
$ws.Range("C6").Value = "https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt"
$ws.Range("D6").Value = "'4.319"
$ws.Range("E6").Value = "'0.03%"
$ws.Range("B7").Value = "FTXToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D7").Value = "'1.524"
$ws.Range("E7").Value = "'-7.24%"
$ws.Range("B8").Value = "MXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx"
$ws.Range("D8").Value = "'0.9227"
$ws.Range("E8").Value = "'0.91%"
$ws.Range("B9").Value = "BTSEToken"
$ws.Range("C9").Value = "https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("D9").Value = "'2.399"
$ws.Range("E9").Value = "'-1.28%"
$ws.Range("B10").Value = "LiechtensteinCryptoassetsExchange"
$ws.Range("C10").Value = "https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("D10").Value = "'0.1166"
$ws.Range("E10").Value = "'-6.35%"
$ws.Range("B11").Value = "WazirX"
$ws.Range("C11").Value = "https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx"
$ws.Range("D11").Value = "'0.1754"
$ws.Range("E11").Value = "'-4.01%"
$ws.Range("B12").Value = "MandalaExchangeToken"
$ws.Range("C12").Value = "https://coinranking.com/coin/lviNIbma2Xuqs+mandalaexchangetoken-mdx"
$ws.Range("D12").Value = "'0.08637"
$ws.Range("E12").Value = "'-6.09%"
$ws.Range("B13").Value = "BitrueCoin"
$ws.Range("C13").Value = "https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("D13").Value = "'0.04181"
$ws.Range("E13").Value = "'-0.99%"
$ws.Range("B14").Value = "BitMartToken"
$ws.Range("C14").Value = "https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("D14").Value = "'0.1053"
$ws.Range("E14").Value = "'0.14%"
$ws.Range("B15").Value = "BitForexToken"
$ws.Range("C15").Value = "https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("D15").Value = "'0.001274"
$ws.Range("E15").Value = "'1.38%"
$ws.Range("B16").Value = "TigerCash"
$ws.Range("C16").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D16").Value = "'0.005883"
$ws.Range("E16").Value = "'1.60%"
$ws.Range("B17").Value = "LEO"
$ws.Range("C17").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D17").Value = "'3.370"
$ws.Range("E17").Value = "'0.72%"
$ws.Range("E18").Value = "'-0.67%"
$ws.Range("D19").Value = "'7.598"
$ws.Range("E19").Value = "'2.25%"
$ws.Range("D20").Value = "'0.1357"
$ws.Range("E20").Value = "'-1.89%"
$ws.Range("D22").Value = "'0.03834"
$ws.Range("E22").Value = "'-4.76%"
$ws.Range("D23").Value = "'0.001285"
$ws.Range("E23").Value = "'1.53%"
$ws.Range("D24").Value = "'0.003622"
$ws.Range("E24").Value = "'-11.50%"
$ws.Range("D25").Value = "'0.0001307"
$ws.Range("E25").Value = "'0.49%"
$ws.Range("D26").Value = "'0.0003740"
$ws.Range("E26").Value = "'-95.02%"
$ws.Range("D38").Value = "'0.02316"
$ws.Range("E38").Value = "'-10.03%"
$ws.Range("D39").Value = "'0.04995"
$ws.Range("E39").Value = "'-6.49%"
$ws.Range("D40").Value = "'0.007714"
$ws.Range("E40").Value = "'-1.72%"
$ws.Range("E41").Value = "'-3.10%"
$ws.Range("E42").Value = "'115.78%"
$ws.Range("D43").Value = "'0.007418"
$ws.Range("E43").Value = "'10.95%"
$ws.Range("D44").Value = "'0.007910"
$ws.Range("E44").Value = "'-1.75%"
$ws.Range("D45").Value = "'0.3165"
$ws.Range("E45").Value = "'3.10%"
$ws.Range("D46").Value = "'0.00006487"
$ws.Range("E46").Value = "'-3.76%"
$ws.Range("E47").Value = "'0.36%"
$ws.Range("E48").Value = "'3.10%"
$ws.Range("D49").Value = "'0.004220"
$ws.Range("E49").Value = "'35.98%"
$ws.Range("D50").Value = "'0.00002110"
$ws.Range("E50").Value = "'0.36%"
$ws.Range("D51").Value = "'0.0002010"
$ws.Range("E51").Value = "'0.36%"
